$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the IP value in C2 (was 127.0.0.1, now 192.168.1.113)
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "192.168.1.113"

# Adjust column C width to fit new content (target stored width = 15)
$ws.Columns("C").ColumnWidth = 14.285714285714286

# Update the active selection to C2
$ws.Range("C2").Select()
